# "added issues to to-do list"
#
# The To-do List - Progress Tracker workbook tracks tasks grouped into
# "Server" (and other) sections. The "Server" section (header row 61,
# items starting row 62) had two still-empty item rows (62 and 63) that
# are now filled in with two new issues:
#   Row 62: Task "Setup Azure account", Owner "FC", 100% complete, dated 6/13/2012
#   Row 63: Task "Install Azure SDK",  Owner "all" (not yet started)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 62 - Setup Azure account, owned by FC, fully complete, dated 2012-06-13
$ws.Range("B62").Value = "Setup Azure account"
$ws.Range("C62").Value = "FC"
$ws.Range("D62").Value = 1
$ws.Range("E62").Value = "6/13/2012"

# Row 63 - Install Azure SDK, owned by everyone ("all"), not started yet
$ws.Range("B63").Value = "Install Azure SDK"
$ws.Range("C63").Value = "all"

# Leave the cursor where the author left off editing the new rows.
$ws.Range("G67").Select()
